$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All changed cells originate from text ("inlineStr") cells in the source.
# Force text number format before assignment so Excel does not auto-convert
# numeric-looking strings (e.g. "0.860", "12.20", "0.0580") into numbers,
# which would silently drop meaningful trailing/leading zeros.

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = '@'
$cell.Value = '35.482.90'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.61%  '
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.896.08'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.69%  '
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.83%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '248.02'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.10%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.691'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.15%  '
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.87%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '44.06'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +8.27%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.353'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.63%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0743'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.59%  '
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.87%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '13.14'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.09%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.168.56'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.81%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.727'
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.14%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.94'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.55%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.905.43'
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.47%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '35.440.09'
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.51%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '73.91'
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.22%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0825'
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.05%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '247.92'
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.65%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '12.88'
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.35%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.98'
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.54%  '
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.75%  '
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.29%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.20'
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -9.83%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '165.83'
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.25%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '8.49'
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.32%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '18.41'
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.89%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.128'
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.52%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.128.41'
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.02%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.78'
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +8.96%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.26'
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.79%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0580'
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.39%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.23'
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.15%  '
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.88%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.860'
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -6.04%  '
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.75%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.58'
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -20.97%  '
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.14%  '
$cell = $ws.Cells.Item(40, 2)
$cell.NumberFormat = '@'
$cell.Value = 'Aave'
$cell = $ws.Cells.Item(40, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '98.01'
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.43%  '
$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = '@'
$cell.Value = 'InjectiveProtocol'
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '17.17'
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.62%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0214'
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.41%  '
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.46%  '
$cell = $ws.Cells.Item(44, 2)
$cell.NumberFormat = '@'
$cell.Value = 'RenderToken'
$cell = $ws.Cells.Item(44, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.38'
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.03%  '
$cell = $ws.Cells.Item(45, 2)
$cell.NumberFormat = '@'
$cell.Value = 'Maker'
$cell = $ws.Cells.Item(45, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.297.41'
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.98%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0799'
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +6.22%  '
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.09%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.75'
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.65%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '12.20'
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.14%  '
$cell = $ws.Cells.Item(50, 2)
$cell.NumberFormat = '@'
$cell.Value = 'MultiversX'
$cell = $ws.Cells.Item(50, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '43.69'
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.02%  '
$cell = $ws.Cells.Item(51, 2)
$cell.NumberFormat = '@'
$cell.Value = 'FraxShare'
$cell = $ws.Cells.Item(51, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.36'
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.43%  '
